# RTM.xlsx update: refresh the newly added High Level Design entries for
# Register, Login and Logout (per Omar's request), and tidy up the
# corresponding Low Level Design cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register -------------------------------------------------------
# New High Level Design function replaces the old placeholder text.
$ws.Range("F6").Value = "High_Design_navigate_to_view_account_page_1,`nHigh_Design_navigate_to_admin_page_1"
$ws.Range("F3").Value = "High_level_navigate_after_register_1"

# --- Login ------------------------------------------------------------
# Low Level Design cell had a stray duplicate reference; clean it up.
$ws.Range("E6").Value = "Low_Design_login_1,`nLow_Design_login_admin_1,`nLow_Design_login_client_2"

# --- Logout -------------------------------------------------------
$ws.Range("F32").Value = "High_Design_Logout_1"

# --- Register tidy up --------------------------------------------------
$ws.Range("E3").Value = "Low_Design_reg_1"

# The Logout Low Level Design cell is no longer needed.
$ws.Range("E32").ClearContents()

# New blank styled row appended below the table (matches the same
# formatting as the last data row).
$ws.Range("F32").Copy($ws.Range("F33"))
$ws.Range("F33").ClearContents()
$ws.Rows.Item(33).RowHeight = 51.75

# Restore the view: scroll down near the bottom of the table and leave
# the selection on the last edited cell.
$ws.Range("F32").Select()
